$d = $word.ActiveDocument

# ------------------------------------------------------------------
# Part 1: "Q-Q plot" -> "histogram" and the expanded sentence about
# the distribution of differences / sample means.
# ------------------------------------------------------------------

# 1a. Swap "Q-Q plot" for "histogram" in the first run of the sentence.
$f1 = $d.Content
$f1.Find.Execute("Q-Q plot", $false, $false, $false, $false, $false, $true, 1, $false, "histogram", 2)

# 1b. Expand "the distribution of sample means is normal:" into the
#     longer sentence that also re-mentions "distribution sample means".
$f2 = $d.Content
$f2.Find.Execute(
    "the distribution of sample means is normal:",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "the distribution of differences is normal in order to be able to claim that the distribution sample means is normal:",
    2)

# At this point Word has coalesced the three original runs ("b. We have...
# to be sure" / " " / "the distribution...normal:") into a single run,
# because they all share the same (default) run formatting and the edit
# touched that shared run-group. We now need to re-split that single run
# back into five runs matching the target layout:
#   [b. We have ... to be sure] [ ] [the distribution ... the] [ ] [distribution sample means is normal:]
# Toggling a character property on/off over a sub-range forces Word to
# split the run at that boundary without altering any text.

# Locate the run that now holds the whole merged sentence.
$fStart = $d.Content
$fStart.Find.Execute("b. We have a small sample", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$runStart = $fStart.Start

# Boundary 1: right after "...to be sure" (start of the first space run).
$fA = $d.Content
$fA.Find.Execute("to be sure", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$boundaryA = $fA.End

# Boundary 2: right after the first space (start of "the distribution...").
$fB = $d.Content
$fB.Find.Execute("the distribution of differences", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$boundaryB = $fB.Start

# Boundary 3: right after "...claim that the" (start of the second space run).
$fC = $d.Content
$fC.Find.Execute("claim that the", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$boundaryC = $fC.End

# Boundary 4: right after the second space (start of "distribution sample means...").
$fD = $d.Content
$fD.Find.Execute("distribution sample means is normal:", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$boundaryD = $fD.Start

# Force splits at each boundary (earliest last, so earlier offsets stay valid).
$sD = $d.Range($runStart, $boundaryD)
$sD.Bold = 1
$sD.Bold = 0

$sC = $d.Range($runStart, $boundaryC)
$sC.Bold = 1
$sC.Bold = 0

$sB = $d.Range($runStart, $boundaryB)
$sB.Bold = 1
$sB.Bold = 0

$sA = $d.Range($runStart, $boundaryA)
$sA.Bold = 1
$sA.Bold = 0

# ------------------------------------------------------------------
# Part 2: "The data appear to be normal" -> "The differences appear to
# be normal"
# ------------------------------------------------------------------
$f3 = $d.Content
$f3.Find.Execute(
    "The data appear to be normal",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "The differences appear to be normal",
    2)
